$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 17
$ws.Cells.Item($row, 1).Value = "r585"
$ws.Cells.Item($row, 2).Value = "timmy"
$ws.Cells.Item($row, 3).Value = "water bath way too hot"
$ws.Cells.Item($row, 4).Value = "2025-10-01 16:11:15"
